# Apply updated Betfair Back/Lay odds to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Australian A-League Men: Macarthur FC vs Auckland FC)
$ws.Range("H2").Value  = 2.28
$ws.Range("J2").Value  = 3.75
$ws.Range("L2").Value  = 1.34
$ws.Range("X2").Value  = 19
$ws.Range("Y2").Value  = 14
$ws.Range("Z2").Value  = 15.5
$ws.Range("AA2").Value = 34
$ws.Range("AB2").Value = 17
$ws.Range("AC2").Value = 8.4
$ws.Range("AD2").Value = 11.5
$ws.Range("AE2").Value = 23
$ws.Range("AF2").Value = 25
$ws.Range("AG2").Value = 14.5
$ws.Range("AH2").Value = 17.5
$ws.Range("AI2").Value = 34
$ws.Range("AK2").Value = 36
$ws.Range("AL2").Value = 42
$ws.Range("AM2").Value = 80
$ws.Range("AN2").Value = 32

# Row 3 (Portuguese Segunda Liga: Benfica B vs Porto B)
$ws.Range("G3").Value = 2.08

# Row 4
$ws.Range("J4").Value = 5.2
$ws.Range("P4").Value = 2.84

# Row 5
$ws.Range("I5").Value = 7.2
$ws.Range("J5").Value = 3.7

# Row 6
$ws.Range("F6").Value = 2.32
$ws.Range("G6").Value = 2.72
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 4.1
$ws.Range("J6").Value = 3
$ws.Range("P6").Value = 1.67

# Row 7
$ws.Range("G7").Value = 1.74
$ws.Range("J7").Value = 3.55
$ws.Range("P7").Value = 1.69

# Row 8
$ws.Range("J8").Value = 2.94
$ws.Range("K8").Value = 5
$ws.Range("P8").Value = 1.65
$ws.Range("Q8").Value = 2

# Row 9
$ws.Range("H9").Value = 3.75
$ws.Range("I9").Value = 5.6
$ws.Range("K9").Value = 5.8
$ws.Range("P9").Value = 1.69
$ws.Range("Q9").Value = 1.93
